$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.876.51'
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = '2.908.98'
$ws.Range('E3').Value = '  -1.79%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.55'
$ws.Range('E5').Value = '  -4.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.83'
$ws.Range('E6').Value = '  -1.91%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('D9').Value = '2.904.06'
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.63'
$ws.Range('E10').Value = '  -8.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.148'
$ws.Range('E11').Value = '  -2.49%  '
$ws.Range('E12').Value = '  -2.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000233'
$ws.Range('E13').Value = '  -3.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.25'
$ws.Range('E14').Value = '  -2.76%  '
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('D16').Value = '3.394.46'
$ws.Range('E16').Value = '  -1.69%  '
$ws.Range('D17').Value = '61.886.57'
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('E18').Value = '  -1.16%  '
$ws.Range('D19').Value = '2.894.18'
$ws.Range('E19').Value = '  -2.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '437.82'
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.32'
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.660'
$ws.Range('E22').Value = '  -1.36%  '
$ws.Range('E23').Value = '  -2.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.46'
$ws.Range('E24').Value = '  -2.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.92'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  -8.04%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  -3.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000105'
$ws.Range('E29').Value = '  +7.45%  '
$ws.Range('E30').Value = '  -2.49%  '
$ws.Range('E31').Value = '  -4.07%  '
$ws.Range('E32').Value = '  -4.50%  '
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  -3.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.959'
$ws.Range('E36').Value = '  -3.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.46'
$ws.Range('E37').Value = '  -3.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.02'
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.96'
$ws.Range('E39').Value = '  -4.28%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.85'
$ws.Range('E40').Value = '  -9.52%  '
$ws.Range('E41').Value = '  -1.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.31'
$ws.Range('E42').Value = '  -2.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.79'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D45').Value = '2.701.47'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '133.95'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '339.95'
$ws.Range('E49').Value = '  -5.78%  '
$ws.Range('E50').Value = '  -1.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.84'
$ws.Range('E51').Value = '  -4.56%  '
